$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 93 - In re Sand Hill Exchange, et al.
$ws.Range("D93").Value = "Settlement"
$ws.Range("E93").Value = "Unregistered Exchange"
$ws.Range("F93").Value = "Civil"
$ws.Range("G93").Value = "N/A"
$ws.Range("H93").Value = "Sand Hill Exchange"
$ws.Range("I93").Value = "N/A"
$ws.Range("J93").Value = 20000
$ws.Range("K93").Value = 1
$ws.Range("L93").Value = 1
$ws.Range("M93").Value = "Washington, D.C."

# Row 94 - In re BTC Trading, Corp. and Ethan Burnside
$ws.Range("D94").Value = "Settlement"
$ws.Range("E94").Value = "Unregistered Exchange"
$ws.Range("F94").Value = "Civil"
$ws.Range("G94").Value = "N/A"
$ws.Range("H94").Value = "BTC Trading Corp."
$ws.Range("I94").Value = "Bitcoin"
$ws.Range("J94").Value = 68000
$ws.Range("K94").Value = 1
$ws.Range("L94").Value = 1
$ws.Range("M94").Value = "New York"

# Row 95 - In re Erik T. Voorhees
$ws.Range("D95").Value = "Settlement"
$ws.Range("E95").Value = "Unregistered Offering"
$ws.Range("F95").Value = "Civil"
$ws.Range("G95").Value = "N/A"
$ws.Range("H95").Value = "FeedZeBirds and SatoshiDICE"
$ws.Range("I95").Value = "Bitcoin"
$ws.Range("J95").Value = 50843.92
$ws.Range("K95").Value = 1
$ws.Range("L95").Value = 0
$ws.Range("M95").Value = "New York"

# Row 96 - SEC v. Shavers (note: column D intentionally left blank, matching source)
$ws.Range("E96").Value = "Unregistered Offering"
$ws.Range("F96").Value = "Civil"
$ws.Range("G96").Value = "N/A"
$ws.Range("H96").Value = "Bitcoin Savings and Trust"
$ws.Range("I96").Value = "Bitcoin"
$ws.Range("J96").Value = 4500000
$ws.Range("K96").Value = 1
$ws.Range("L96").Value = 1
$ws.Range("M96").Value = "Washington, D.C."

# Update the view: scroll/selection moved to F97 with no frozen top-left cell
$ws.Range("F97").Select()
